$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'recovery compression pants'
$ws.Cells.Item(2, 1).Value = 'knee volleyball'
$ws.Cells.Item(3, 1).Value = 'leggings black'
$ws.Cells.Item(4, 1).Value = 'mens volleyball pads'
$ws.Cells.Item(5, 1).Value = 'knee arthritis compression'
$ws.Cells.Item(6, 1).Value = 'youth small baseball pants'
$ws.Cells.Item(7, 1).Value = 'mens black pants'
$ws.Cells.Item(8, 1).Value = 'knee high leggings'
$ws.Cells.Item(9, 1).Value = 'knee pads arthritis'
$ws.Cells.Item(10, 1).Value = 'basketball knee pads nike pro'
$ws.Cells.Item(11, 1).Value = 'knee compression tights'
$ws.Cells.Item(12, 1).Value = 'knee pad sleeve youth'
$ws.Cells.Item(13, 1).Value = 'thermal compression pants men winter'
$ws.Cells.Item(14, 1).Value = 'honey comb knee pad'
$ws.Cells.Item(15, 1).Value = 'combat pants with knee pads'
$ws.Cells.Item(16, 1).Value = 'deportivo para hombres'
$ws.Cells.Item(17, 1).Value = 'compression spats'
$ws.Cells.Item(18, 1).Value = 'red mens compression pants'
$ws.Cells.Item(19, 1).Value = 'green basketball knee pads'
$ws.Cells.Item(20, 1).Value = 'red basketball leggings'
$ws.Cells.Item(21, 1).Value = 'underarmour mens tights'
$ws.Cells.Item(22, 1).Value = 'three quarter compression pants men'
$ws.Cells.Item(23, 1).Value = 'base layer pants men cold weather'
$ws.Cells.Item(24, 1).Value = 'tesla base layer pants'
$ws.Cells.Item(25, 1).Value = 'winter gym leggings'
$ws.Cells.Item(26, 1).Value = 'thick spandex pants men'
$ws.Cells.Item(27, 1).Value = 'knee pad leggings for basketball'
$ws.Cells.Item(28, 1).Value = 'youth knee pads for basketball'
$ws.Cells.Item(29, 1).Value = 'girls basketball knee pad leggings'
$ws.Cells.Item(30, 1).Value = 'knee pads compression'
$ws.Cells.Item(31, 1).Value = 'knee pad tights basketball'
$ws.Cells.Item(32, 1).Value = 'weightlifting knee pads'
$ws.Cells.Item(33, 1).Value = 'compression leggings youth'
$ws.Cells.Item(34, 1).Value = 'hex pad knee'
$ws.Cells.Item(35, 1).Value = 'knee pad soccer'
$ws.Cells.Item(36, 1).Value = 'knee pad for running'
$ws.Cells.Item(37, 1).Value = 'basketballs black'
$ws.Cells.Item(38, 1).Value = 'boys soccer tights'
$ws.Cells.Item(39, 1).Value = 'compression men leggings'
$ws.Cells.Item(40, 1).Value = 'youth compression knee pads'
$ws.Cells.Item(41, 1).Value = 'leggings with knee'
$ws.Cells.Item(42, 1).Value = 'running compression pants men'
$ws.Cells.Item(43, 1).Value = 'knee protector volleyball'
$ws.Cells.Item(44, 1).Value = 'volleyball knee pads for men'
$ws.Cells.Item(45, 1).Value = 'knee guard basketball'
$ws.Cells.Item(46, 1).Value = 'black baseball pants boys'
$ws.Cells.Item(47, 1).Value = 'compression pants boys'
$ws.Cells.Item(48, 1).Value = 'knee pad adult'
$ws.Cells.Item(49, 1).Value = 'little boy leggings'
$ws.Cells.Item(50, 1).Value = 'capri legging pack'
$ws.Cells.Item(51, 1).Value = 'men knee pads'
$ws.Cells.Item(52, 1).Value = 'basketball compression'
$ws.Cells.Item(53, 1).Value = 'boys compression tight'
$ws.Cells.Item(54, 1).Value = 'baseball pants knee high'
$ws.Cells.Item(55, 1).Value = 'small volleyball knee pads'
$ws.Cells.Item(56, 1).Value = 'knee pads volleyball'
$ws.Cells.Item(57, 1).Value = 'youth volleyball'
$ws.Cells.Item(58, 1).Value = 'knee protector football'
$ws.Cells.Item(59, 1).Value = 'knee pads protector'
$ws.Cells.Item(60, 1).Value = 'snowboarding hip pads'
$ws.Cells.Item(61, 1).Value = 'knee pad running'
$ws.Cells.Item(62, 1).Value = 'men pads'
$ws.Cells.Item(63, 1).Value = 'cycling compression tights'
$ws.Cells.Item(64, 1).Value = 'capri leggings pack'
$ws.Cells.Item(65, 1).Value = 'pants soccer'
$ws.Cells.Item(66, 1).Value = 'running compression leggings'
$ws.Cells.Item(67, 1).Value = 'boys black baseball pants'
$ws.Cells.Item(68, 1).Value = 'leggings for boys'
$ws.Cells.Item(69, 1).Value = 'youth girls compression pants'
$ws.Cells.Item(70, 1).Value = 'athletic pants soccer'
$ws.Cells.Item(71, 1).Value = 'knee guards basketball'
$ws.Cells.Item(72, 1).Value = 'compression pants for hockey'
$ws.Cells.Item(73, 1).Value = 'wrestling gear for men'
$ws.Cells.Item(74, 1).Value = 'adult medium baseball pants'
$ws.Cells.Item(75, 1).Value = 'knee pads for arthritis'
$ws.Cells.Item(76, 1).Value = 'baseball boy pants'
$ws.Cells.Item(77, 1).Value = 'knee protector pad'
$ws.Cells.Item(78, 1).Value = 'knee protector soccer'
$ws.Cells.Item(79, 1).Value = 'men sports compression'
$ws.Cells.Item(80, 1).Value = 'basketball pants boys'
$ws.Cells.Item(81, 1).Value = 'baseball pants'
$ws.Cells.Item(82, 1).Value = 'knee protector pads'
$ws.Cells.Item(83, 1).Value = 'knee pad for sports'
$ws.Cells.Item(84, 1).Value = 'sports knee protectors'
$ws.Cells.Item(85, 1).Value = 'youth knee pads volleyball'
$ws.Cells.Item(86, 1).Value = 'bjj pants'
$ws.Cells.Item(87, 1).Value = 'knee pads cycling'
$ws.Cells.Item(88, 1).Value = 'boys youth basketball'
$ws.Cells.Item(89, 1).Value = 'mens athletic pants tall'
$ws.Cells.Item(90, 1).Value = 'legging pack'
$ws.Cells.Item(91, 1).Value = 'girl volleyball knee pads'
$ws.Cells.Item(92, 1).Value = 'running tights'
$ws.Cells.Item(93, 1).Value = 'youth baseball pants'
$ws.Cells.Item(94, 1).Value = 'boys sports pants'
$ws.Cells.Item(95, 1).Value = 'compression for men'
$ws.Cells.Item(96, 1).Value = 'basketballs youth size'
$ws.Cells.Item(97, 1).Value = 'knee pads xl'
$ws.Cells.Item(98, 1).Value = 'compression knee men'
$ws.Cells.Item(99, 1).Value = 'knee guards for adults'
$ws.Cells.Item(100, 1).Value = 'football knee pad'
